$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New s_val data (filtered to exclude save games), row by row B:G
$data = @{
    2 = @(1.445647641019636, 1.626987699542094, 0.1496068669990043, 0.5333859586016987, 0, 3.755628166162433)
    3 = @(0.6545652718822623, 1.626987699542094, 18.71679738969934, 0.5333859586016987, 1, 21.53173631972539)
    4 = @(1.445647641019636, 1.626987699542094, 0.1496068669990043, 0.5333859586016987, 1, 3.755628166162433)
    5 = @(0.6545652718822623, 0.3048912486333797, 0.7210945179870265, 0.5333859586016987, 1, 2.213936997104367)
    6 = @(0.6545652718822623, 0.3048912486333797, 0.7210945179870265, 0.5333859586016987, 1, 2.213936997104367)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 5).Value = $vals[3]
    $ws.Cells.Item($row, 6).Value = $vals[4]
    $ws.Cells.Item($row, 7).Value = $vals[5]
}

$wb.Save()
